# Generate Report for Handoff
# A new handoff (c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md /
# ffff6c0ace36-61ad-4947-a5d9-19faa014a1d7.md) superseded the previous one
# (2e39905d-8c4d-4586-9f6a-548e46376d5a.md / 7801fcc7-0d93-4ab8-953e-b281839d361e.md).
# This updates the Overview / zh-cn / de-de sheets to reflect the new
# "Ready for handoff" status, new handoff file/datetime, and clears the
# stale "Latest Target File" / "Latest Handback File" columns + resets the
# handback datetime, since nothing has been handed back yet for this handoff.

$wb = $excel.ActiveWorkbook

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/92d5a1f02dbd3e2ee1936909132cff6efe05fd8c/e2e/"
$zhHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/88fc9541de2b46b30e173b884f759669e94173cc/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deHandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ee8b1fa82ac8d78e2ae3d60dab9e97eaaf854d5b/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

$newFile1 = "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.md"
$newFile2 = "ffff6c0ace36-61ad-4947-a5d9-19faa014a1d7.md"
$zhXlf = "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.zh-cn.xlf"
$deXlf = "c85dd3f2-ec34-4d7e-980f-d33a2a27cfba.eb7edf2a86468b4bd614ade89c8221dda5c35aab.de-de.xlf"

$status = "Ready for handoff"
$overviewDate = "2016-03-22 21:11:02"
$zhHandoffDate = "2016-03-22 21:10:58"
$nullDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = $status
$wsOverview.Range("C2").Value = $status
$wsOverview.Range("D2").Value = $overviewDate

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("D3").Value = $overviewDate

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), ($mdBase + $newFile1), [System.Type]::Missing, [System.Type]::Missing, $newFile1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), ($mdBase + $newFile2), [System.Type]::Missing, [System.Type]::Missing, $newFile2)

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFile1
$wsZh.Range("C2").Value = $status
$wsZh.Range("D2").Value = $zhXlf
$wsZh.Range("E2").Value = $zhHandoffDate
$wsZh.Range("F2").Clear()
$wsZh.Range("G2").Clear()
$wsZh.Range("H2").Value = $nullDate

$wsZh.Range("A3").Value = $newFile2
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = $zhXlf
$wsZh.Range("E3").Value = $zhHandoffDate
$wsZh.Range("F3").Clear()
$wsZh.Range("G3").Clear()
$wsZh.Range("H3").Value = $nullDate

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($mdBase + $newFile1), [System.Type]::Missing, [System.Type]::Missing, $newFile1)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), ($zhHandoffBase + $zhXlf), [System.Type]::Missing, [System.Type]::Missing, $zhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($mdBase + $newFile2), [System.Type]::Missing, [System.Type]::Missing, $newFile2)
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), ($zhHandoffBase + $zhXlf), [System.Type]::Missing, [System.Type]::Missing, $zhXlf)

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFile1
$wsDe.Range("C2").Value = $status
$wsDe.Range("D2").Value = $deXlf
$wsDe.Range("E2").Value = $overviewDate
$wsDe.Range("F2").Clear()
$wsDe.Range("G2").Clear()
$wsDe.Range("H2").Value = $nullDate

$wsDe.Range("A3").Value = $newFile2
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = $deXlf
$wsDe.Range("E3").Value = $overviewDate
$wsDe.Range("F3").Clear()
$wsDe.Range("G3").Clear()
$wsDe.Range("H3").Value = $nullDate

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($mdBase + $newFile1), [System.Type]::Missing, [System.Type]::Missing, $newFile1)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), ($deHandoffBase + $deXlf), [System.Type]::Missing, [System.Type]::Missing, $deXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($mdBase + $newFile2), [System.Type]::Missing, [System.Type]::Missing, $newFile2)
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), ($deHandoffBase + $deXlf), [System.Type]::Missing, [System.Type]::Missing, $deXlf)

Write-Host "Generated report for handoff."
